# Refresh monthly budget figures on the "manoj.o" sheet to match the
# latest uploaded source data, and update the view state (scroll
# position + selection) to where the author left off editing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manoj.o")
$ws.Activate()

# Row 2
$ws.Cells.Item(2, 3).Value = 35000
$ws.Cells.Item(2, 4).Value = 35000
$ws.Cells.Item(2, 5).Value = 35000
$ws.Cells.Item(2, 6).Value = 35000
$ws.Cells.Item(2, 7).Value = 35000
$ws.Cells.Item(2, 8).Value = 35000
$ws.Cells.Item(2, 9).Value = 35000
$ws.Cells.Item(2, 10).Value = 35000
$ws.Cells.Item(2, 11).Value = 35000
$ws.Cells.Item(2, 12).Value = 35000
$ws.Cells.Item(2, 13).Value = 35000
$ws.Cells.Item(2, 14).Value = 35000

# Row 3
$ws.Cells.Item(3, 11).Value = 0

# Row 4
$ws.Cells.Item(4, 3).Value = 100
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 500
$ws.Cells.Item(4, 10).Value = 200

# Row 5
$ws.Cells.Item(5, 3).Value = 1000
$ws.Cells.Item(5, 4).Value = 1000
$ws.Cells.Item(5, 5).Value = 1000
$ws.Cells.Item(5, 6).Value = 1000
$ws.Cells.Item(5, 7).Value = 1000
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 12).Value = 1000
$ws.Cells.Item(5, 13).Value = 1000
$ws.Cells.Item(5, 14).Value = 1000

# Row 6
$ws.Cells.Item(6, 12).Value = 4500

# Row 8
$ws.Cells.Item(8, 4).Value = 20000
$ws.Cells.Item(8, 9).Value = 20000

# Row 10
$ws.Cells.Item(10, 5).Value = 5000

# Row 11
$ws.Cells.Item(11, 3).Value = 1000
$ws.Cells.Item(11, 4).Value = 1000
$ws.Cells.Item(11, 5).Value = 1000
$ws.Cells.Item(11, 6).Value = 1000
$ws.Cells.Item(11, 7).Value = 1000
$ws.Cells.Item(11, 8).Value = 1000
$ws.Cells.Item(11, 9).Value = 1000
$ws.Cells.Item(11, 10).Value = 1000
$ws.Cells.Item(11, 11).Value = 1000
$ws.Cells.Item(11, 12).Value = 1000
$ws.Cells.Item(11, 13).Value = 1000
$ws.Cells.Item(11, 14).Value = 1000

# Row 13
$ws.Cells.Item(13, 3).Value = 1400
$ws.Cells.Item(13, 4).Value = 1200
$ws.Cells.Item(13, 5).Value = 2000
$ws.Cells.Item(13, 6).Value = 3000
$ws.Cells.Item(13, 7).Value = 3000
$ws.Cells.Item(13, 8).Value = 2500
$ws.Cells.Item(13, 9).Value = 2000
$ws.Cells.Item(13, 10).Value = 1200
$ws.Cells.Item(13, 11).Value = 1300
$ws.Cells.Item(13, 12).Value = 1250
$ws.Cells.Item(13, 13).Value = 1150
$ws.Cells.Item(13, 14).Value = 1210

# Row 14
$ws.Cells.Item(14, 3).Value = 100
$ws.Cells.Item(14, 4).Value = 100
$ws.Cells.Item(14, 5).Value = 100
$ws.Cells.Item(14, 6).Value = 100
$ws.Cells.Item(14, 7).Value = 100
$ws.Cells.Item(14, 8).Value = 100
$ws.Cells.Item(14, 9).Value = 100
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 100
$ws.Cells.Item(14, 12).Value = 100
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(14, 14).Value = 100

# Row 15
$ws.Cells.Item(15, 5).Value = 2000
$ws.Cells.Item(15, 8).Value = 2000
$ws.Cells.Item(15, 11).Value = 2000
$ws.Cells.Item(15, 14).Value = 2000

# Row 16
$ws.Cells.Item(16, 3).Value = 200
$ws.Cells.Item(16, 5).Value = 200
$ws.Cells.Item(16, 7).Value = 200
$ws.Cells.Item(16, 9).Value = 200
$ws.Cells.Item(16, 11).Value = 200
$ws.Cells.Item(16, 13).Value = 200

# Row 17
$ws.Cells.Item(17, 8).Value = 6000
$ws.Cells.Item(17, 14).Value = 6000

# Row 18
$ws.Cells.Item(18, 3).Value = 2000
$ws.Cells.Item(18, 4).Value = 2000
$ws.Cells.Item(18, 5).Value = 2000
$ws.Cells.Item(18, 6).Value = 2000
$ws.Cells.Item(18, 7).Value = 2000
$ws.Cells.Item(18, 8).Value = 2000
$ws.Cells.Item(18, 9).Value = 2000
$ws.Cells.Item(18, 10).Value = 2000
$ws.Cells.Item(18, 11).Value = 2000
$ws.Cells.Item(18, 12).Value = 2000
$ws.Cells.Item(18, 13).Value = 2000
$ws.Cells.Item(18, 14).Value = 2000

# Row 19
$ws.Cells.Item(19, 3).Value = 1500
$ws.Cells.Item(19, 4).Value = 1500
$ws.Cells.Item(19, 5).Value = 1500
$ws.Cells.Item(19, 6).Value = 1500
$ws.Cells.Item(19, 7).Value = 1500
$ws.Cells.Item(19, 8).Value = 1500
$ws.Cells.Item(19, 9).Value = 1500
$ws.Cells.Item(19, 10).Value = 1500
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1500
$ws.Cells.Item(19, 14).Value = 1500

# Row 20
$ws.Cells.Item(20, 5).Value = 1350
$ws.Cells.Item(20, 8).Value = 1350
$ws.Cells.Item(20, 11).Value = 1350
$ws.Cells.Item(20, 14).Value = 1350

# Row 21
$ws.Cells.Item(21, 3).Value = 800
$ws.Cells.Item(21, 4).Value = 800
$ws.Cells.Item(21, 5).Value = 800
$ws.Cells.Item(21, 6).Value = 800
$ws.Cells.Item(21, 7).Value = 800
$ws.Cells.Item(21, 8).Value = 800
$ws.Cells.Item(21, 9).Value = 800
$ws.Cells.Item(21, 10).Value = 800
$ws.Cells.Item(21, 11).Value = 800
$ws.Cells.Item(21, 12).Value = 800
$ws.Cells.Item(21, 13).Value = 800
$ws.Cells.Item(21, 14).Value = 800

# Row 22
$ws.Cells.Item(22, 3).Value = 100
$ws.Cells.Item(22, 4).Value = 150
$ws.Cells.Item(22, 5).Value = 100
$ws.Cells.Item(22, 6).Value = 200
$ws.Cells.Item(22, 7).Value = 100
$ws.Cells.Item(22, 8).Value = 50
$ws.Cells.Item(22, 9).Value = 100
$ws.Cells.Item(22, 10).Value = 150
$ws.Cells.Item(22, 11).Value = 125
$ws.Cells.Item(22, 12).Value = 130
$ws.Cells.Item(22, 13).Value = 140
$ws.Cells.Item(22, 14).Value = 110

# Row 23
$ws.Cells.Item(23, 4).Value = 10000

# Row 24
$ws.Cells.Item(24, 3).Value = 10000
$ws.Cells.Item(24, 4).Value = 10000
$ws.Cells.Item(24, 5).Value = 10000
$ws.Cells.Item(24, 6).Value = 10000
$ws.Cells.Item(24, 7).Value = 10000
$ws.Cells.Item(24, 8).Value = 10000
$ws.Cells.Item(24, 9).Value = 11000
$ws.Cells.Item(24, 10).Value = 11000
$ws.Cells.Item(24, 11).Value = 11000
$ws.Cells.Item(24, 12).Value = 11000
$ws.Cells.Item(24, 13).Value = 11000
$ws.Cells.Item(24, 14).Value = 11000

# Row 25
$ws.Cells.Item(25, 3).Value = 100
$ws.Cells.Item(25, 4).Value = 200
$ws.Cells.Item(25, 5).Value = 100
$ws.Cells.Item(25, 6).Value = 100
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 50
$ws.Cells.Item(25, 11).Value = 20
$ws.Cells.Item(25, 12).Value = 30
$ws.Cells.Item(25, 13).Value = 100

# Row 26
$ws.Cells.Item(26, 4).Value = 50

# Row 27
$ws.Cells.Item(27, 4).Value = 50000

# Row 28
$ws.Cells.Item(28, 3).Value = 2000
$ws.Cells.Item(28, 4).Value = 2000
$ws.Cells.Item(28, 5).Value = 2000
$ws.Cells.Item(28, 6).Value = 2000
$ws.Cells.Item(28, 7).Value = 2000
$ws.Cells.Item(28, 8).Value = 2000
$ws.Cells.Item(28, 9).Value = 2000
$ws.Cells.Item(28, 10).Value = 2000
$ws.Cells.Item(28, 11).Value = 2000
$ws.Cells.Item(28, 12).Value = 2000
$ws.Cells.Item(28, 13).Value = 2000
$ws.Cells.Item(28, 14).Value = 2000

# Row 29
$ws.Cells.Item(29, 3).Value = 150
$ws.Cells.Item(29, 4).Value = 150
$ws.Cells.Item(29, 5).Value = 150
$ws.Cells.Item(29, 6).Value = 150
$ws.Cells.Item(29, 7).Value = 150
$ws.Cells.Item(29, 8).Value = 150
$ws.Cells.Item(29, 9).Value = 150
$ws.Cells.Item(29, 10).Value = 150
$ws.Cells.Item(29, 11).Value = 150
$ws.Cells.Item(29, 12).Value = 150
$ws.Cells.Item(29, 13).Value = 150
$ws.Cells.Item(29, 14).Value = 150

# Row 30
$ws.Cells.Item(30, 3).Value = 500
$ws.Cells.Item(30, 4).Value = 2000
$ws.Cells.Item(30, 5).Value = 500
$ws.Cells.Item(30, 6).Value = 500
$ws.Cells.Item(30, 7).Value = 500
$ws.Cells.Item(30, 8).Value = 500
$ws.Cells.Item(30, 9).Value = 500
$ws.Cells.Item(30, 10).Value = 2000
$ws.Cells.Item(30, 11).Value = 500
$ws.Cells.Item(30, 12).Value = 5000
$ws.Cells.Item(30, 13).Value = 500
$ws.Cells.Item(30, 14).Value = 500

# Row 31
$ws.Cells.Item(31, 4).Value = 60
$ws.Cells.Item(31, 6).Value = 60
$ws.Cells.Item(31, 8).Value = 100
$ws.Cells.Item(31, 10).Value = 60
$ws.Cells.Item(31, 12).Value = 100
$ws.Cells.Item(31, 14).Value = 60

# Row 32
$ws.Cells.Item(32, 3).Value = 200
$ws.Cells.Item(32, 4).Value = 200
$ws.Cells.Item(32, 5).Value = 200
$ws.Cells.Item(32, 6).Value = 200
$ws.Cells.Item(32, 7).Value = 200
$ws.Cells.Item(32, 8).Value = 200
$ws.Cells.Item(32, 9).Value = 200
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 200
$ws.Cells.Item(32, 12).Value = 200
$ws.Cells.Item(32, 13).Value = 200
$ws.Cells.Item(32, 14).Value = 200

# Row 33
$ws.Cells.Item(33, 3).Value = 300
$ws.Cells.Item(33, 4).Value = 300
$ws.Cells.Item(33, 5).Value = 500
$ws.Cells.Item(33, 6).Value = 300
$ws.Cells.Item(33, 7).Value = 300
$ws.Cells.Item(33, 8).Value = 500
$ws.Cells.Item(33, 9).Value = 300
$ws.Cells.Item(33, 10).Value = 300
$ws.Cells.Item(33, 11).Value = 500
$ws.Cells.Item(33, 12).Value = 300
$ws.Cells.Item(33, 13).Value = 500
$ws.Cells.Item(33, 14).Value = 1000

# Row 34
$ws.Cells.Item(34, 3).Value = 200
$ws.Cells.Item(34, 4).Value = 200
$ws.Cells.Item(34, 5).Value = 200
$ws.Cells.Item(34, 6).Value = 200
$ws.Cells.Item(34, 7).Value = 200
$ws.Cells.Item(34, 8).Value = 200
$ws.Cells.Item(34, 9).Value = 200
$ws.Cells.Item(34, 10).Value = 200
$ws.Cells.Item(34, 11).Value = 200
$ws.Cells.Item(34, 12).Value = 200
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 200

# Row 35
$ws.Cells.Item(35, 4).Value = 300
$ws.Cells.Item(35, 6).Value = 300
$ws.Cells.Item(35, 7).Value = 400
$ws.Cells.Item(35, 9).Value = 250
$ws.Cells.Item(35, 11).Value = 300

# Row 36
$ws.Cells.Item(36, 8).Value = 2000
$ws.Cells.Item(36, 9).Value = 5000
$ws.Cells.Item(36, 12).Value = 1000

# Row 37
$ws.Cells.Item(37, 3).Value = 200
$ws.Cells.Item(37, 8).Value = 100

# Row 38
$ws.Cells.Item(38, 10).Value = 3000
$ws.Cells.Item(38, 14).Value = 2000

# Row 39
$ws.Cells.Item(39, 3).Value = 1000
$ws.Cells.Item(39, 4).Value = 1000
$ws.Cells.Item(39, 5).Value = 1000
$ws.Cells.Item(39, 6).Value = 1000
$ws.Cells.Item(39, 7).Value = 1000
$ws.Cells.Item(39, 8).Value = 1000
$ws.Cells.Item(39, 9).Value = 1000
$ws.Cells.Item(39, 10).Value = 1000
$ws.Cells.Item(39, 11).Value = 1000
$ws.Cells.Item(39, 12).Value = 1000
$ws.Cells.Item(39, 13).Value = 1000
$ws.Cells.Item(39, 14).Value = 1000

# Scroll the window so column B is the left-most visible column and
# row 16 is the top-most visible row (mirrors topLeftCell "B16"),
# then select the final totals row as the active selection.
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D39:N39").Select()
